$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price figures stored as text (values can contain
# thousands separators like "26.483.57", which are not valid numbers).
# Force the whole data range to Text format first so plain decimal
# looking values (e.g. "219.56") are not silently reinterpreted as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.483.57"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "1.678.41"
$ws.Range("E3").Value = "  +2.50%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "219.56"
$ws.Range("E5").Value = "  +2.74%  "
$ws.Range("E6").Value = "  +2.13%  "
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "0.2715"
$ws.Range("E8").Value = "  +4.67%  "
$ws.Range("E9").Value = "  +2.15%  "
$ws.Range("D10").Value = "21.95"
$ws.Range("E10").Value = "  +6.64%  "
$ws.Range("D11").Value = "0.07800"
$ws.Range("E11").Value = "  +1.59%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.683.20"
$ws.Range("E12").Value = "  +2.27%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.521"
$ws.Range("E13").Value = "  +2.79%  "
$ws.Range("D14").Value = "0.5604"
$ws.Range("E14").Value = "  +1.84%  "
$ws.Range("D15").Value = "0.0₅8352"
$ws.Range("E15").Value = "  +1.83%  "
$ws.Range("D16").Value = "65.83"
$ws.Range("E16").Value = "  +1.55%  "
$ws.Range("D17").Value = "26.523.40"
$ws.Range("E17").Value = "  +1.77%  "
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("E19").Value = "  +2.97%  "
$ws.Range("D20").Value = "193.56"
$ws.Range("E20").Value = "  +3.05%  "
$ws.Range("E21").Value = "  +1.42%  "
$ws.Range("D22").Value = "6.329"
$ws.Range("E22").Value = "  +2.98%  "
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E24").Value = "  -2.08%  "
$ws.Range("D25").Value = "0.1282"
$ws.Range("E25").Value = "  +6.05%  "
$ws.Range("D26").Value = "7.423"
$ws.Range("E26").Value = "  +0.47%  "
$ws.Range("E27").Value = "  +3.48%  "
$ws.Range("D28").Value = "1.444"
$ws.Range("E28").Value = "  +3.66%  "
$ws.Range("D29").Value = "0.06296"
$ws.Range("E29").Value = "  +5.95%  "
$ws.Range("D30").Value = "1.288"
$ws.Range("E30").Value = "  +2.73%  "
$ws.Range("D31").Value = "3.611"
$ws.Range("E31").Value = "  +5.47%  "
$ws.Range("D32").Value = "3.468"
$ws.Range("E32").Value = "  +2.20%  "
$ws.Range("D33").Value = "1.705"
$ws.Range("E33").Value = "  +3.97%  "
$ws.Range("E34").Value = "  +3.38%  "
$ws.Range("D35").Value = "0.6154"
$ws.Range("E35").Value = "  +9.45%  "
$ws.Range("D36").Value = "2.420"
$ws.Range("E36").Value = "  +1.06%  "
$ws.Range("D37").Value = "2.785"
$ws.Range("E37").Value = "  +0.87%  "
$ws.Range("B38").Value = "FraxShare"
$ws.Range("C38").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D38").Value = "6.161"
$ws.Range("E38").Value = "  +8.23%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01636"
$ws.Range("E39").Value = "  +1.52%  "
$ws.Range("D40").Value = "1.092.28"
$ws.Range("E40").Value = "  +6.14%  "
$ws.Range("D41").Value = "0.8670"
$ws.Range("E41").Value = "  +2.06%  "
$ws.Range("E42").Value = "  -0.12%  "
$ws.Range("D43").Value = "100.72"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").Value = "1.823.97"
$ws.Range("E44").Value = "  +2.07%  "
$ws.Range("E45").Value = "  +2.78%  "
$ws.Range("D46").Value = "58.96"
$ws.Range("E46").Value = "  +5.98%  "
$ws.Range("D47").Value = "8.194"
$ws.Range("E47").Value = "  +1.76%  "
$ws.Range("D48").Value = "1.002"
$ws.Range("E48").Value = "  +0.33%  "
$ws.Range("E49").Value = "  +1.32%  "
$ws.Range("D50").Value = "6.055"
$ws.Range("E50").Value = "  +2.70%  "
$ws.Range("D51").Value = "1.474"
$ws.Range("E51").Value = "  +6.91%  "
